# Update Leve profit/price figures across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ /
# LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns H-N) per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 928.5714
$ws.Range("I18").Value = 833.3333
$ws.Range("J18").Value = 1500
$ws.Range("K18").Value = 833.3333
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = -549.3333
$ws.Range("N18").Value = -2068
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H51").Value = 67322.53
$ws.Range("J51").Value = 90066.164
$ws.Range("L51").Value = 90066.164
$ws.Range("N51").Value = -91034.164
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H106").Value = 37054396
$ws.Range("I106").Value = 40015504
$ws.Range("K106").Value = 40015504
$ws.Range("M106").Value = -40014873
$ws.Range("H132").Value = 2314.92
$ws.Range("I132").Value = 1221.5454
$ws.Range("K132").Value = 3664.6362
$ws.Range("M132").Value = -1134.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4968.7144
$ws.Range("I32").Value = 3353.3928
$ws.Range("K32").Value = 3353.3928
$ws.Range("M32").Value = -3066.3928
$ws.Range("H61").Value = 3892.2666
$ws.Range("I61").Value = 3941.7144
$ws.Range("J61").Value = 3200
$ws.Range("K61").Value = 3941.7144
$ws.Range("L61").Value = 3200
$ws.Range("M61").Value = -3729.7144
$ws.Range("N61").Value = -3624
$ws.Range("H74").Value = 4145
$ws.Range("I74").Value = 3955.5557
$ws.Range("K74").Value = 3955.5557
$ws.Range("M74").Value = -3081.5557
$ws.Range("H77").Value = 4145
$ws.Range("I77").Value = 3955.5557
$ws.Range("K77").Value = 19777.7785
$ws.Range("M77").Value = -15409.7785
$ws.Range("H110").Value = 3586323.2
$ws.Range("I110").Value = 5557011
$ws.Range("J110").Value = 3254.5454
$ws.Range("K110").Value = 5557011
$ws.Range("L110").Value = 3254.5454
$ws.Range("M110").Value = -5554966
$ws.Range("N110").Value = -7344.5454
$ws.Range("H132").Value = 3948.15
$ws.Range("I132").Value = 4178.909
$ws.Range("K132").Value = 12536.727
$ws.Range("M132").Value = -10006.727
$ws.Range("H136").Value = 3892.2666
$ws.Range("I136").Value = 3941.7144
$ws.Range("J136").Value = 3200
$ws.Range("K136").Value = 11825.1432
$ws.Range("L136").Value = 9600
$ws.Range("M136").Value = -9275.143199999999
$ws.Range("N136").Value = -14700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 11479
$ws.Range("I36").Value = 2958
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 2958
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -2424
$ws.Range("N36").Value = -21068
$ws.Range("H37").Value = 613
$ws.Range("J37").Value = 1000
$ws.Range("L37").Value = 1000
$ws.Range("N37").Value = -1274

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 10171.667
$ws.Range("I3").Value = 7853.8
$ws.Range("K3").Value = 7853.8
$ws.Range("M3").Value = -7740.8
$ws.Range("H15").Value = 5103.3335
$ws.Range("I15").Value = 14420
$ws.Range("J15").Value = 445
$ws.Range("K15").Value = 14420
$ws.Range("L15").Value = 445
$ws.Range("M15").Value = -14250
$ws.Range("N15").Value = -785
$ws.Range("H31").Value = 3297.0645
$ws.Range("I31").Value = 2314.9565
$ws.Range("K31").Value = 2314.9565
$ws.Range("M31").Value = -2019.9565
$ws.Range("H34").Value = 3297.0645
$ws.Range("I34").Value = 2314.9565
$ws.Range("K34").Value = 2314.9565
$ws.Range("M34").Value = -2112.9565
$ws.Range("H122").Value = 1695.7
$ws.Range("I122").Value = 995
$ws.Range("K122").Value = 2985
$ws.Range("M122").Value = -535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 4479.3
$ws.Range("J97").Value = 5099.125
$ws.Range("L97").Value = 15297.375
$ws.Range("N97").Value = -16289.375
$ws.Range("H107").Value = 200621
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 200621
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 601863
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -605703
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()
$ws.Range("H131").Value = 1400.3677
$ws.Range("J131").Value = 1472.0161
$ws.Range("L131").Value = 4416.0483
$ws.Range("N131").Value = -14496.0483

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 43222.25
$ws.Range("J46").Value = 30000
$ws.Range("L46").Value = 30000
$ws.Range("N46").Value = -30312
$ws.Range("H70").Value = 7430.4443
$ws.Range("I70").Value = 7146.5
$ws.Range("K70").Value = 7146.5
$ws.Range("M70").Value = -6876.5
$ws.Range("H73").Value = 7430.4443
$ws.Range("I73").Value = 7146.5
$ws.Range("K73").Value = 7146.5
$ws.Range("M73").Value = -6210.5
$ws.Range("H80").Value = 3564.9443
$ws.Range("I80").Value = 2864.0833
$ws.Range("J80").Value = 4966.6665
$ws.Range("K80").Value = 2864.0833
$ws.Range("L80").Value = 4966.6665
$ws.Range("M80").Value = -1866.0833
$ws.Range("N80").Value = -6962.6665
$ws.Range("H83").Value = 3564.9443
$ws.Range("I83").Value = 2864.0833
$ws.Range("J83").Value = 4966.6665
$ws.Range("K83").Value = 14320.4165
$ws.Range("L83").Value = 24833.3325
$ws.Range("M83").Value = -9328.416499999999
$ws.Range("N83").Value = -34817.3325
$ws.Range("H97").Value = 1828.8
$ws.Range("I97").Value = 1071.4117
$ws.Range("J97").Value = 3438.25
$ws.Range("K97").Value = 1071.4117
$ws.Range("L97").Value = 3438.25
$ws.Range("M97").Value = -575.4117000000001
$ws.Range("N97").Value = -4430.25
$ws.Range("H107").Value = 1149.8334
$ws.Range("I107").Value = 1100.25
$ws.Range("K107").Value = 1100.25
$ws.Range("M107").Value = 819.75
$ws.Range("H132").Value = 1623.7142
$ws.Range("I132").Value = 1623.7142
$ws.Range("K132").Value = 4871.142599999999
$ws.Range("M132").Value = -2341.142599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 15011
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 15011
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 15011
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -15237
$ws.Range("H40").Value = 1195.6
$ws.Range("I40").Value = 1195.6
$ws.Range("K40").Value = 1195.6
$ws.Range("M40").Value = -1059.6
$ws.Range("H41").Value = 33999
$ws.Range("I41").Value = 33999
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 33999
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -33561
$ws.Range("N41").ClearContents()
$ws.Range("H47").Value = 28999
$ws.Range("I47").Value = 28999
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 28999
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -28509
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 28999
$ws.Range("I52").Value = 28999
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 28999
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -28766
$ws.Range("N52").ClearContents()
$ws.Range("H93").Value = 971.9655
$ws.Range("I93").Value = 788.4286
$ws.Range("K93").Value = 788.4286
$ws.Range("M93").Value = 459.5714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1832.4706
$ws.Range("I100").Value = 1537
$ws.Range("K100").Value = 3074
$ws.Range("M100").Value = -2533
$ws.Range("H113").Value = 2366.5
$ws.Range("I113").Value = 839.8
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 2519.4
$ws.Range("L113").Value = 30000
$ws.Range("M113").Value = -349.3999999999996
$ws.Range("N113").Value = -34340
$ws.Range("H126").Value = 1439.4
$ws.Range("I126").Value = 1439.4
$ws.Range("K126").Value = 4318.200000000001
$ws.Range("M126").Value = -1848.200000000001
